$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 13 (STT=9): task name changed
$ws.Range("B13").Value = "Thiết kế Form Information"

# Rows 10-14 (STT 6-10): fill in the "actual end date" column (H) with the text "20/10/2018"
$ws.Range("H10").Value = "20/10/2018"
$ws.Range("H11").Value = "20/10/2018"
$ws.Range("H12").Value = "20/10/2018"
$ws.Range("H13").Value = "20/10/2018"
$ws.Range("H14").Value = "20/10/2018"
